# Stage 보상 api 추가
# 버그 확인 필요

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

# Rename field "WorldId" -> "WorldNum" in the Packet model definition (A4)
$ws.Range("A4").Value = "WorldNum"

# Update the active selection to A5 (as recorded in the saved view state)
$ws.Range("A5").Select()
